# Apply the "include product name" edit to the evaluation workbook.
#
# Summary of the change:
#  - Each product worksheet's header cell (C2) is updated from a short
#    "Product: X" label to the full product name ("<full product name>: X").
#  - Each worksheet's selection moves from A6:XFD6 to C2:G2 (C3:G3 on the
#    Price sheet), i.e. highlighting the header row that was just edited.
#  - Worksheet tab names are renamed from the generic "Product-X" scheme
#    to more descriptive category-based names.
#  - The workbook re-opens with the first sheet active (instead of the
#    last "Edge" sheet), and the stray topLeftCell scroll position on the
#    Pockets sheet is cleared.

$wb = $excel.ActiveWorkbook

# --- Rename worksheet tabs -------------------------------------------------
$wb.Worksheets.Item(1).Name = "BikeAccessory-Grips"
$wb.Worksheets.Item(2).Name = "DivingSet-Mask"
$wb.Worksheets.Item(3).Name = "TensionMeter-Spokes"
$wb.Worksheets.Item(4).Name = "DotSight-Battery"
$wb.Worksheets.Item(5).Name = "BottleCage-Price"
$wb.Worksheets.Item(6).Name = "Backpack-Pockets"
$wb.Worksheets.Item(7).Name = "Knife-Edge"

# --- Update header cell (C2) on every sheet with the full product name ----
$wb.Worksheets.Item(1).Range("C2").Value = "Sunlite MX 1 Foam Comfort Bicycle Grips, Black: Grips"
$wb.Worksheets.Item(2).Range("C2").Value = "U.S. Divers Men's Lux LX Mask with Purge, Pivot Fins and Phoenix LX Snorkel Combo Set: Mask"
$wb.Worksheets.Item(3).Range("C2").Value = "Park Tool TM-1 Spoke Tension Meter: Spokes"
$wb.Worksheets.Item(4).Range("C2").Value = "Leapers Golden Image 38mm Red/Green Dot Sight, Integral Weaver Mount: Battery"
$wb.Worksheets.Item(5).Range("C2").Value = "Ibera Bicycle Lightweight Aluminum Water Bottle Cage: Price"
$wb.Worksheets.Item(6).Range("C2").Value = "5.11 Rush 72 Back Pack: Feature"
$wb.Worksheets.Item(7).Range("C2").Value = "Survivor HK-106320 Outdoor Fixed Blade Knife 7 Overall WITH FIRE STARTER: Edge"

# --- Update selections: header block C2:G2 (C3:G3 on the Price sheet) -----
# (Selecting a fresh range below also clears the stray topLeftCell="A3"
# scroll position that was saved on the Pockets sheet.)
# Activate every other sheet first so that only the first sheet ends up
# marked as the active tab once we activate it last below.
$wb.Worksheets.Item(2).Activate()
$wb.Worksheets.Item(2).Range("C2:G2").Select() | Out-Null

$wb.Worksheets.Item(3).Activate()
$wb.Worksheets.Item(3).Range("C2:G2").Select() | Out-Null

$wb.Worksheets.Item(4).Activate()
$wb.Worksheets.Item(4).Range("C2:G2").Select() | Out-Null

$wb.Worksheets.Item(5).Activate()
$wb.Worksheets.Item(5).Range("C3:G3").Select() | Out-Null

$wb.Worksheets.Item(6).Activate()
$wb.Worksheets.Item(6).Range("C2:G2").Select() | Out-Null

$wb.Worksheets.Item(7).Activate()
$wb.Worksheets.Item(7).Range("C2:G2").Select() | Out-Null

$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("C2:G2").Select() | Out-Null

# --- Restore the workbook window chrome recorded in the saved view -------
$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 0
$win.Width = 25600
$win.Height = 16000
